$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$mapping = @{
    "gpt-5.1" = "Model C"
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro" = "Model B"
    "kimi-k2" = "Model D"
}

for ($row = 2; $row -le 49; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($mapping.ContainsKey($current)) {
        $cell.Value = $mapping[$current]
    }
}
